$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (matches source data which is inline text, not numbers)
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '63.725.29'
$ws.Range("E2").Value = '  -4.00%  '
$ws.Range("D3").Value = '3.111.07'
$ws.Range("E3").Value = '  -4.92%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '607.14'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").Value = '143.83'
$ws.Range("E6").Value = '  -9.48%  '
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").Value = '3.109.66'
$ws.Range("E8").Value = '  -4.88%  '
$ws.Range("D9").Value = '0.516'
$ws.Range("E9").Value = '  -4.82%  '
$ws.Range("D10").Value = '0.149'
$ws.Range("E10").Value = '  -7.47%  '
$ws.Range("E11").Value = '  -8.84%  '
$ws.Range("E12").Value = '  -5.78%  '
$ws.Range("D13").Value = '0.0000249'
$ws.Range("E13").Value = '  -8.42%  '
$ws.Range("D14").Value = '34.95'
$ws.Range("E14").Value = '  -10.10%  '
$ws.Range("D15").Value = '3.637.20'
$ws.Range("E15").Value = '  -4.65%  '
$ws.Range("E16").Value = '  +1.53%  '
$ws.Range("D17").Value = '63.786.13'
$ws.Range("E17").Value = '  -4.00%  '
$ws.Range("D18").Value = '3.124.35'
$ws.Range("E18").Value = '  -4.76%  '
$ws.Range("D19").Value = '6.76'
$ws.Range("E19").Value = '  -8.25%  '
$ws.Range("D20").Value = '475.39'
$ws.Range("E20").Value = '  -5.53%  '
$ws.Range("D21").Value = '14.54'
$ws.Range("E21").Value = '  -5.85%  '
$ws.Range("D22").Value = '0.705'
$ws.Range("E22").Value = '  -6.49%  '
$ws.Range("D23").Value = '7.65'
$ws.Range("E23").Value = '  -5.80%  '
$ws.Range("D24").Value = '13.48'
$ws.Range("E24").Value = '  -8.55%  '
$ws.Range("D25").Value = '83.90'
$ws.Range("E25").Value = '  -3.08%  '
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("D27").Value = '2.77'
$ws.Range("E27").Value = '  -8.87%  '
$ws.Range("D28").Value = '8.36'
$ws.Range("E28").Value = '  -9.07%  '
$ws.Range("D29").Value = '2.09'
$ws.Range("E29").Value = '  -12.51%  '
$ws.Range("D30").Value = '6.74'
$ws.Range("E30").Value = '  -4.42%  '
$ws.Range("E31").Value = '  -17.30%  '
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("D33").Value = '2.67'
$ws.Range("E33").Value = '  -6.53%  '
$ws.Range("D34").Value = '26.10'
$ws.Range("E34").Value = '  -6.60%  '
$ws.Range("D35").Value = '1.11'
$ws.Range("E35").Value = '  -3.08%  '
$ws.Range("D36").Value = '5.88'
$ws.Range("E36").Value = '  -9.21%  '
$ws.Range("D37").Value = '52.78'
$ws.Range("E37").Value = '  -5.21%  '
$ws.Range("D38").Value = '0.0₃0740'
$ws.Range("E38").Value = '  -7.09%  '
$ws.Range("D39").Value = '454.15'
$ws.Range("E39").Value = '  -9.69%  '
$ws.Range("D40").Value = '2.90'
$ws.Range("E40").Value = '  -17.50%  '
$ws.Range("D41").Value = '0.0389'
$ws.Range("E41").Value = '  -8.74%  '
$ws.Range("D42").Value = '0.117'
$ws.Range("E42").Value = '  -9.62%  '
$ws.Range("D43").Value = '8.27'
$ws.Range("E43").Value = '  -6.04%  '
$ws.Range("D44").Value = '2.841.27'
$ws.Range("E44").Value = '  -5.52%  '
$ws.Range("D45").Value = '0.263'
$ws.Range("E45").Value = '  -10.13%  '
$ws.Range("D46").Value = '2.25'
$ws.Range("E46").Value = '  -13.43%  '
$ws.Range("D47").Value = '2.41'
$ws.Range("E47").Value = '  -2.93%  '
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("D49").Value = '25.88'
$ws.Range("E49").Value = '  -10.60%  '
$ws.Range("D50").Value = '0.112'
$ws.Range("E50").Value = '  -5.55%  '
$ws.Range("D51").Value = '118.06'
$ws.Range("E51").Value = '  -2.81%  '

# Restore default style on column D so only the cell *content* changed (avoid leaving a stray number-format override)
$dRange.Style = "Normal"
